# Add a new homework/data row (row 82) for 2020-02-13, as per commit
# "Added homework 13Feb 2020".
#
# The sheet is a simple stock-price table with columns:
#   A timestamp, B date, C id, D name, E open, F high, G low, H close, I vol
#
# Columns B (date) and C (id) are stored as *text* in the workbook (e.g.
# "2020-02-13" and "0212"), not as real Excel dates / numbers. Assigning a
# plain string to Range.Value lets Excel's smart text recognition kick in
# and silently convert "2020-02-13" into a date serial number and "0212"
# into the number 212 (dropping the leading zero). To keep these as literal
# text - matching every other row in the column - we build them via a
# formula that yields a text result, then convert that formula to a static
# value with Copy + PasteSpecial(xlPasteValues). That preserves the text
# type without Excel re-interpreting the literal, and (unlike forcing a
# "Text" NumberFormat or a leading apostrophe) it does not leave behind an
# extra cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 82

$ws.Cells.Item($row, 1).Value = 1581552000

$ws.Cells.Item($row, 2).Formula = '="2020-02-13"'
$ws.Cells.Item($row, 2).Copy()
$ws.Cells.Item($row, 2).PasteSpecial(-4163)

$ws.Cells.Item($row, 3).Formula = '="0212"'
$ws.Cells.Item($row, 3).Copy()
$ws.Cells.Item($row, 3).PasteSpecial(-4163)

$ws.Cells.Item($row, 4).Value = "SDS"
$ws.Cells.Item($row, 5).Value = 0.225
$ws.Cells.Item($row, 6).Value = 0.235
$ws.Cells.Item($row, 7).Value = 0.22
$ws.Cells.Item($row, 8).Value = 0.23
$ws.Cells.Item($row, 9).Value = 499300
